# Franchise added and register page re-designed, register request needs
# testing, including fav franchise setting.
#
# This script:
#   1. Inserts a new "Franchises" worksheet right before "SQL Commands"
#      and populates it with FranchiseID / Location2022 / Nickname2022
#      data (mirrors the 2022 season rows already on the "Teams" sheet)
#      plus a helper column that builds the SQL INSERT statement.
#   2. Adds a "FranchiseID" header to the "Teams" sheet (column K) so the
#      new Franchises table can eventually be linked back to Teams.
#   3. Moves the active-tab/selection state from "UserScores" to "Teams"
#      (re-designed register page references Teams/Franchises).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Franchises" sheet, inserted just before "SQL Commands"
# ---------------------------------------------------------------------
$sqlSheet = $wb.Worksheets.Item("SQL Commands")
$franchises = $wb.Worksheets.Add($sqlSheet)
$franchises.Name = "Franchises"

$franchises.Range("A1").Value = "FranchiseID"
$franchises.Range("B1").Value = "Location2022"
$franchises.Range("C1").Value = "Nickname2022"
$franchises.Range("E1").Value = "INSERT INTO Franchises (FranchiseID, Location2022, Nickname2022) VALUES"
$franchises.Range("A1:C1").Font.Bold = $true
$franchises.Range("E1").Font.Bold = $true

$franchiseData = @(
  @(1, "Baltimore", "Ravens"),
  @(2, "Cincinnati", "Bengals"),
  @(3, "Cleveland", "Browns"),
  @(4, "Pittsburgh", "Steelers"),
  @(5, "Buffalo", "Bills"),
  @(6, "Miami", "Dolphins"),
  @(7, "New England", "Patriots"),
  @(8, "New York", "Jets"),
  @(9, "Houston", "Texans"),
  @(10, "Indianapolis", "Colts"),
  @(11, "Jacksonville", "Jaguars"),
  @(12, "Tennessee", "Titans"),
  @(13, "Denver", "Broncos"),
  @(14, "Kansas City", "Chiefs"),
  @(15, "Las Vegas", "Raiders"),
  @(16, "Los Angeles", "Chargers"),
  @(17, "Chicago", "Bears"),
  @(18, "Detroit", "Lions"),
  @(19, "Green Bay", "Packers"),
  @(20, "Minnesota", "Vikings"),
  @(21, "Dallas", "Cowboys"),
  @(22, "New York", "Giants"),
  @(23, "Philadelphia", "Eagles"),
  @(24, "Washington", "Commanders"),
  @(25, "Atlanta", "Falcons"),
  @(26, "Carolina", "Panthers"),
  @(27, "New Orleans", "Saints"),
  @(28, "Tampa Bay", "Buccaneers"),
  @(29, "Arizona", "Cardinals"),
  @(30, "Los Angeles", "Rams"),
  @(31, "San Francisco", "49ers"),
  @(32, "Seattle", "Seahawks")
)

$r = 2
foreach ($row in $franchiseData) {
  $franchises.Range("A$r").Value = $row[0]
  $franchises.Range("B$r").Value = $row[1]
  $franchises.Range("C$r").Value = $row[2]
  $r++
}

$franchises.Range("E2").Formula = '="("&A2&",''"&B2&"'',''"&C2&"''),"'
$franchises.Range("E3:E33").Formula = '="("&A3&",''"&B3&"'',''"&C3&"''),"'

$franchises.Columns.Item(1).ColumnWidth = 10.67
$franchises.Columns.Item(2).ColumnWidth = 12.33
$franchises.Columns.Item(3).ColumnWidth = 13.67

# ---------------------------------------------------------------------
# 2. "Teams" sheet gains a FranchiseID header (column K), matching the
#    other bold headers in row 1
# ---------------------------------------------------------------------
$teams = $wb.Worksheets.Item("Teams")
$teams.Range("K1").Value = "FranchiseID"
$teams.Range("K1").Font.Bold = $true

# ---------------------------------------------------------------------
# 3. Selection / active-tab bookkeeping: move away from "UserScores"
#    and onto "Teams" (re-designed register page). Each sheet keeps its
#    own last selection; only the final Activate()/Select() decides
#    which tab ends up "active" for the whole workbook.
# ---------------------------------------------------------------------
$userScores = $wb.Worksheets.Item("UserScores")
$userScores.Range("D1:D111").Select() | Out-Null

$franchises.Activate()
$franchises.Range("K22").Select() | Out-Null

$teams.Activate()
$teams.Range("M9").Select() | Out-Null
